# "Latex chapter 1 and revision"
#
# The sheet previously had its table starting at row 3 (rows 1-2 were
# blank), leaving the used range as A3:D31. This edit removes the two
# leading blank rows so the table now starts at row 1 (A1:D29), shifting
# every data row up by two - the row contents themselves are unchanged.
# It also updates the active selection left behind on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two empty leading rows; Excel shifts rows 3-31 up to 1-29,
# preserving every cell's value/style.
[void]$ws.Rows("1:2").Delete()

# Leave the selection where the edit was made (column D, header+first row).
[void]$ws.Range("D3:D4").Select()
